# "bold, italic and bold&Italic done"
#
# Replace the single auto-numbered list textbox on slide 1 with a new
# plain textbox containing one paragraph of four runs that progressively
# apply italic / bold / bold+italic formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the existing "Text 0" shape (numbered list placeholder-ish text box).
$old = $s.Shapes.Item(1)
$old.Delete()

# EMU -> point conversion (914400 EMU per inch, 72 points per inch).
$emuPerPt = 914400 / 72

$left   = 1448873 / $emuPerPt
$top    = 457199  / $emuPerPt
$width  = 6561786 / $emuPerPt
$height = 369332  / $emuPerPt

# New textbox at the target position/size.
$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)

$tf = $tb.TextFrame
$tr = $tf.TextRange

$run1 = "This word is "
$run2 = "italic "
$run3 = "bold "
$run4 = "finally"

$tr.Text = $run1 + $run2 + $run3 + $run4

# Italicize "italic "
$start2 = $run1.Length + 1
$tr.Characters($start2, $run2.Length).Font.Italic = $true

# Bold "bold "
$start3 = $run1.Length + $run2.Length + 1
$tr.Characters($start3, $run3.Length).Font.Bold = $true

# Bold + italic "finally"
$start4 = $run1.Length + $run2.Length + $run3.Length + 1
$run4Range = $tr.Characters($start4, $run4.Length)
$run4Range.Font.Bold = $true
$run4Range.Font.Italic = $true

# Word-wrap + shape autofit-to-text, matching the textbox's bodyPr.
$tf.WordWrap = -1
$tf.AutoSize = 1

# No shape fill (<a:noFill/>).
$tb.Fill.Visible = 0
